$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 74, shifting existing rows 74-172 down to 75-173.
$ws.Rows(74).Insert()

# Populate the newly inserted row 74 with the new record.
$ws.Cells.Item(74, 1).Value  = 5
$ws.Cells.Item(74, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value  = "Maule"
$ws.Cells.Item(74, 4).Value  = 44579
$ws.Cells.Item(74, 5).Value  = 7
$ws.Cells.Item(74, 6).Value  = 100112024
$ws.Cells.Item(74, 7).Value  = "Choclo"
$ws.Cells.Item(74, 8).Value  = "Choclero"
$ws.Cells.Item(74, 9).Value  = "Primera"
$ws.Cells.Item(74, 10).Value = 45000
$ws.Cells.Item(74, 11).Value = 200
$ws.Cells.Item(74, 12).Value = 200
$ws.Cells.Item(74, 13).Value = 200
$ws.Cells.Item(74, 14).Value = "`$/unidad"
$ws.Cells.Item(74, 15).Value = "Región del Maule"
$ws.Cells.Item(74, 16).Value = 200
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"

Write-Output "Row 74 inserted and populated."
